# lab1_report_example.docx edit ("add hmetis result to lab1_report"):
#  - fill in the Student ID
#  - fill in the Name (Chinese + English)
#  - fill in the shmetis benchmark results (size / time) in both result
#    tables, replacing the "xxx" placeholders

$d = $word.ActiveDocument

# Replace the first occurrence of $oldText inside $range with $newText,
# leaving the rest of the range's content untouched.
function Replace-InRange($range, $oldText, $newText) {
    $full = $range.Text
    $idx = $full.IndexOf($oldText)
    if ($idx -lt 0) {
        throw "text '$oldText' not found in range '$full'"
    }
    $base = $range.Start
    $target = $d.Range($base + $idx, $base + $idx + $oldText.Length)
    $target.Text = $newText
}

# ---------------------------------------------------------------------
# 1) "Student_ID: XXXXXXX" -> "Student_ID: 0710006"
# ---------------------------------------------------------------------
Replace-InRange $d.Paragraphs.Item(1).Range "XXXXXXX" "0710006"

# ---------------------------------------------------------------------
# 2) "Name: XXX" -> "Name: 盧可瑜 Ke-Yu Lu"
# ---------------------------------------------------------------------
Replace-InRange $d.Paragraphs.Item(2).Range "XXX" "盧可瑜 Ke-Yu Lu"

# ---------------------------------------------------------------------
# 3) First results table (Case1 / Case2 / Case3 columns) - shmetis row:
#    fill in size/runtime pairs for each case.
# ---------------------------------------------------------------------
$table1 = $d.Tables.Item(1)
$shmetisRow1 = 3

Replace-InRange $table1.Cell($shmetisRow1, 2).Range "xxx" "1"
Replace-InRange $table1.Cell($shmetisRow1, 3).Range "xxx" "0.001"
Replace-InRange $table1.Cell($shmetisRow1, 4).Range "xxx" "33"
Replace-InRange $table1.Cell($shmetisRow1, 5).Range "xxx" "0.310"
Replace-InRange $table1.Cell($shmetisRow1, 6).Range "xxx" "199"
Replace-InRange $table1.Cell($shmetisRow1, 7).Range "xxx" "1.602"

# ---------------------------------------------------------------------
# 4) Second results table (Case4 column) - shmetis row.
# ---------------------------------------------------------------------
$table2 = $d.Tables.Item(2)
$shmetisRow2 = 3

Replace-InRange $table2.Cell($shmetisRow2, 2).Range "xxx" "18"
Replace-InRange $table2.Cell($shmetisRow2, 3).Range "xxx" "4.804"
